# Updated symbol list with refreshed prices and 1h volume percentages for cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "E2" "1.89%"
Set-TextValue "D3" "27.32"
Set-TextValue "E3" "1.42%"
Set-TextValue "D4" "4.710"
Set-TextValue "E4" "-1.15%"
Set-TextValue "D5" "0.06090"
Set-TextValue "E5" "3.00%"
Set-TextValue "D6" "6.679"
Set-TextValue "E6" "1.05%"
Set-TextValue "D7" "0.8467"
Set-TextValue "E7" "-0.49%"
Set-TextValue "D8" "0.9289"
Set-TextValue "E8" "0.78%"
Set-TextValue "D9" "0.1406"
Set-TextValue "E9" "2.00%"
Set-TextValue "D10" "0.04773"
Set-TextValue "E10" "14.90%"
Set-TextValue "D11" "0.07103"
Set-TextValue "E11" "1.47%"
Set-TextValue "D12" "0.03091"
Set-TextValue "E12" "1.36%"
Set-TextValue "D13" "0.09065"
Set-TextValue "E13" "-0.39%"
Set-TextValue "D14" "0.001542"
Set-TextValue "E14" "0.08%"
Set-TextValue "D15" "0.0006098"
Set-TextValue "E15" "0.66%"
Set-TextValue "D16" "0.006143"
Set-TextValue "E16" "1.87%"
Set-TextValue "D17" "3.448"
Set-TextValue "E17" "-0.69%"
Set-TextValue "E18" "-0.56%"
Set-TextValue "D19" "2.163"
Set-TextValue "E19" "-0.64%"
Set-TextValue "E20" "2.25%"
Set-TextValue "E21" "0.22%"
Set-TextValue "D22" "4.086"
Set-TextValue "E22" "4.72%"
Set-TextValue "D23" "0.04233"
Set-TextValue "E23" "-0.49%"
Set-TextValue "D25" "0.003796"
Set-TextValue "E25" "-11.22%"
Set-TextValue "D26" "0.0001201"
Set-TextValue "E26" "0.08%"
Set-TextValue "D40" "0.03877"
Set-TextValue "E40" "2.44%"
Set-TextValue "E41" "1.36%"
Set-TextValue "D42" "0.004081"
Set-TextValue "E42" "-34.80%"
Set-TextValue "D43" "0.01627"
Set-TextValue "E43" "15.34%"
Set-TextValue "D45" "0.00005137"
Set-TextValue "E45" "-3.90%"
Set-TextValue "E46" "0.06%"
Set-TextValue "D47" "0.1355"
Set-TextValue "E47" "-39.93%"
Set-TextValue "E49" "0.06%"
Set-TextValue "E50" "0.06%"
